$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("product_life")
$ws.Columns.Item(3).Delete()
